# "Add Uncertainty to Model"
# Apply an updated (reduced) estimate for the container height and
# diameter, reflecting added uncertainty in the dimension assumptions.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newDimension = 1.4137154999999999

$ws.Range("container_height").Value = $newDimension
$ws.Range("container_diameter").Value = $newDimension
